$wb = $excel.ActiveWorkbook

# --- Controllers sheet ---
$wsControllers = $wb.Worksheets.Item("Controllers")
$wsControllers.Range("D70").Value = 1
$wsControllers.Range("D71").Value = 1
$wsControllers.Range("D94").Value = 0.16

# --- Daos sheet ---
$wsDaos = $wb.Worksheets.Item("Daos")
$wsDaos.Range("C21").Value = 0.11
$wsDaos.Range("C73").Value = 1
$wsDaos.Range("C77").Value = 1

# --- Vistas sheet ---
$wsVistas = $wb.Worksheets.Item("Vistas")
$wsVistas.Range("C33").Value = 0.1
$wsVistas.Range("C36").Value = 0.95

# --- Avance sheet ---
$wsAvance = $wb.Worksheets.Item("Avance")
$wsAvance.Range("D11").Value = 5

# --- View / selection state ---
# Final view order matters: last activated/selected sheet becomes the
# workbook's active tab. Target state: Daos is the active tab, with
# Controllers, Vistas and Avance holding their own last-selected cell.

$wsControllers.Activate()
$wsControllers.Range("D67").Select()

$wsVistas.Activate()
$wsVistas.Range("C84").Select()

$wsAvance.Activate()
$wsAvance.Range("D11").Select()

$wsDaos.Activate()
$wsDaos.Range("C74").Select()
